# Update the nested for-loop example (slide 11) and its R console
# output (slide 12): the loop's print() call now labels the outer
# index "k" instead of "i", and the printed output lines follow suit.

$p = $ppt.ActivePresentation

# Slide 11: "Nested for Loop Example" - the print(paste(...)) statement
# lives in paragraph 3 of the body placeholder (shape 2).
$codeSlide = $p.Slides.Item(11)
$codeBody = $codeSlide.Shapes.Item(2).TextFrame.TextRange
$codeBody.Paragraphs(3).Runs(1).Text = '    print(paste("k = " , i, "; j = ", j))'

# Slide 12: "Output" - each console line is its own paragraph in the
# body placeholder (shape 2); update the "i = " label to "k = " while
# keeping the printed numbers untouched.
$outputSlide = $p.Slides.Item(12)
$outputBody = $outputSlide.Shapes.Item(2).TextFrame.TextRange

$outputLines = @(
    '[1] "k =  1 ; j =  1"',
    '[1] "k =  1 ; j =  2"',
    '[1] "k =  2 ; j =  1"',
    '[1] "k =  2 ; j =  2"',
    '[1] "k =  3 ; j =  1"',
    '[1] "k =  3 ; j =  2"',
    '[1] "k =  4 ; j =  1"',
    '[1] "k =  4 ; j =  2"'
)

for ($i = 1; $i -le $outputLines.Length; $i++) {
    $outputBody.Paragraphs($i).Runs(1).Text = $outputLines[$i - 1]
}
